$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.343.11'
$ws.Range('E2').Value = '  +1.86%  '

$ws.Range('D3').Value = '3.930.92'
$ws.Range('E3').Value = '  -0.14%  '

$ws.Range('E4').Value = '  +0.26%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '487.25'

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '148.60'
$ws.Range('E6').Value = '  +2.26%  '

$ws.Range('E7').Value = '  +1.20%  '

$ws.Range('E8').Value = '  +0.03%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.736'
$ws.Range('E9').Value = '  +0.81%  '

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.170'
$ws.Range('E10').Value = '  +3.89%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0000354'
$ws.Range('E11').Value = '  +4.71%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '43.08'
$ws.Range('E12').Value = '  -0.31%  '

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '10.69'
$ws.Range('E13').Value = '  +3.48%  '

$ws.Range('D14').Value = '4.562.44'
$ws.Range('E14').Value = '  +0.31%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '14.55'
$ws.Range('E15').Value = '  -3.72%  '

$ws.Range('D16').Value = '3.924.15'
$ws.Range('E16').Value = '  +0.03%  '

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.136'
$ws.Range('E17').Value = '  -0.58%  '

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '20.02'
$ws.Range('E18').Value = '  +1.11%  '

$ws.Range('E19').Value = '  -1.48%  '

$ws.Range('D20').Value = '68.479.13'
$ws.Range('E20').Value = '  +1.85%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '442.68'
$ws.Range('E21').Value = '  +2.93%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '3.53'
$ws.Range('E22').Value = '  +4.60%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '15.08'
$ws.Range('E23').Value = '  +3.63%  '

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '88.51'
$ws.Range('E24').Value = '  +1.69%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '11.39'
$ws.Range('E25').Value = '  +18.06%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '11.45'
$ws.Range('E26').Value = '  +12.41%  '

$ws.Range('E27').Value = '  +2.61%  '

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '39.00'
$ws.Range('E28').Value = '  +1.35%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '5.83'
$ws.Range('E29').Value = '  +1.68%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '722.48'
$ws.Range('E30').Value = '  -0.75%  '

$ws.Range('E31').Value = '  +1.47%  '

$ws.Range('E32').Value = '  -1.06%  '

$ws.Range('E33').Value = '  +4.28%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '6.29'
$ws.Range('E34').Value = '  +17.19%  '

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '42.42'
$ws.Range('E35').Value = '  -0.67%  '

$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0898'
$ws.Range('E36').Value = '  +15.02%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '61.33'
$ws.Range('E37').Value = '  +5.91%  '

$ws.Range('E38').Value = '  -1.63%  '

$ws.Range('E39').Value = '  +19.97%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '3.07'
$ws.Range('E40').Value = '  +18.41%  '

$ws.Range('E41').Value = '  +0.16%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '3.23'
$ws.Range('E42').Value = '  +6.30%  '

$ws.Range('E43').Value = '  +1.72%  '

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '2.94'
$ws.Range('E44').Value = '  +4.83%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0370'
$ws.Range('E45').Value = '  +44.21%  '

$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.143'
$ws.Range('E46').Value = '  +1.49%  '

$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.00'
$ws.Range('E47').Value = '  +0.15%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '3.28'
$ws.Range('E48').Value = '  +3.56%  '

$ws.Range('B49').Value = 'LidoDAOToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '3.43'
$ws.Range('E49').Value = '  +0.30%  '

$ws.Range('E50').Value = '  -1.48%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '145.59'
$ws.Range('E51').Value = '  -0.57%  '
